$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "CaixaDeTexto 19" (Informações/Comportamento bullets) ---
$shInfo = $s.Shapes.Item(6)
$trInfo = $shInfo.TextFrame.TextRange

# Remove the first two bullets ("Gênero masculino;" and "Média de 30 anos;")
# leaving: "Formado na área de TI;", "Trabalha como suporte técnico;",
# "Trabalha no segmento de call center;"
$null = $trInfo.Characters(1, 36).Delete()

# Add a new bullet, followed by a new (blank) trailing paragraph
$null = $trInfo.InsertAfter("`rPossui uma boa habilidade analítica;`r")

# The shape auto-fits its text; restore its original height (unchanged by
# the source edit) now that the extra paragraphs have been added.
$shInfo.Height = 149.08181102362204

# --- Shape "CaixaDeTexto 21" (Dores e Necessidades bullets) ---
$shDores = $s.Shapes.Item(7)
$trDores = $shDores.TextFrame.TextRange

# Append two new bullets
$null = $trDores.InsertAfter("`rOutros funcionários têm costume de solicitar um atendimento sem efetuar o registro da solicitação;`rProcesso de atendimento definido.")

# Resize the text box to fit the new content
$shDores.Width = 816.8155905511811
$shDores.Height = 178.1631
